$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values to reflect the repulled/recalculated data
$ws.Range("F2").Value = -2
$ws.Range("F8").Value = -1
$ws.Range("F16").Value = -2
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = -5
$ws.Range("F22").Value = 7
$ws.Range("F24").Value = -1
